$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130; this shifts the existing rows
# 130..184 down to 131..185 (matching the diff where every record from
# the old row 130 onward moves one row down).
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new weekly price record.
$ws.Cells.Item(130, 1).Value = 4
$ws.Cells.Item(130, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(130, 3).Value = "Los Lagos"
$ws.Cells.Item(130, 4).Value = 44726
$ws.Cells.Item(130, 5).Value = 10
$ws.Cells.Item(130, 6).Value = 100112009
$ws.Cells.Item(130, 7).Value = "Acelga"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 80
$ws.Cells.Item(130, 11).Value = 12000
$ws.Cells.Item(130, 12).Value = 12000
$ws.Cells.Item(130, 13).Value = 12000
$ws.Cells.Item(130, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(130, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(130, 16).Value = 1000
$ws.Cells.Item(130, 17).Value = 12
$ws.Cells.Item(130, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# the rest of the column (style index 2 in the original workbook).
$ws.Cells.Item(130, 4).NumberFormat = $ws.Cells.Item(131, 4).NumberFormat
